# productListings.xlsx — restore the correct product rows after an earlier
# mis-shuffle: row 14/15, 43/44, 106/107, 113/114 and 117/118 each had their
# productName/productURL/productID swapped with their neighbour, and the
# productPrice in row 13/14 was likewise swapped.
#
# $q is a literal leading apostrophe: prefixing a numeric-looking value with
# it forces Excel to store/keep the cell as TEXT (matching the workbook's
# original inlineStr/text cells) instead of silently re-typing it as a
# number, which would both reformat "2,299.99"-style prices and truncate the
# 19-digit productID values to double precision.
$q = "'"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- rows 13-15 ---------------------------------------------------------
$ws.Range("D13").Value = $q + "129.99"

$ws.Range("B14").Value = "HP ENVY 5055 All-in-One Printer"
$ws.Range("C14").Value = "https://store.hp.com//us/en/pdp/hp-envy-5055-all-in-one-printer"
$ws.Range("D14").Value = $q + "2,299.99"
$ws.Range("E14").Value = $q + "3074457345618778821"

$ws.Range("B15").Value = "HP PageWide Enterprise Color MFP 586f"
$ws.Range("C15").Value = "https://store.hp.com//us/en/pdp/hp-pagewide-enterprise-color-mfp-586f"
$ws.Range("E15").Value = $q + "1282656"

# --- rows 43-44 ----------------------------------------------------------
$ws.Range("B43").Value = "HP OfficeJet 3830 All-in-One Printer"
$ws.Range("C43").Value = "https://store.hp.com//us/en/pdp/hp-officejet-3830-all-in-one-printer"
$ws.Range("D43").Value = $q + "99.99"
$ws.Range("E43").Value = $q + "1030654"

$ws.Range("B44").Value = "HP PageWide Pro 577dw Multifunction Printer"
$ws.Range("C44").Value = "https://store.hp.com//us/en/pdp/hp-pagewide-pro-577dw-multifunction-printer"
$ws.Range("D44").Value = $q + "899.99"
$ws.Range("E44").Value = $q + "1243169"

# --- rows 106-107 ----------------------------------------------------------
$ws.Range("B106").Value = "HP DesignJet T830 24-in Multifunction Printer"
$ws.Range("C106").Value = "https://store.hp.com//us/en/pdp/hp-designjet-t830-24-in-multifunction-printer"
$ws.Range("D106").Value = $q + "2,695.99"
$ws.Range("E106").Value = $q + "3074457345618628324"

$ws.Range("B107").Value = "HP LaserJet Enterprise M607n"
$ws.Range("C107").Value = "https://store.hp.com//us/en/pdp/hp-laserjet-enterprise-m607n"
$ws.Range("D107").Value = $q + "749.99"
$ws.Range("E107").Value = $q + "1695183"

# --- rows 113-114 ----------------------------------------------------------
$ws.Range("B113").Value = "HP LaserJet Pro M203dw Printer"
$ws.Range("C113").Value = "https://store.hp.com//us/en/pdp/hp-laserjet-pro-m203dw-printer"
$ws.Range("D113").Value = $q + "169.99"
$ws.Range("E113").Value = $q + "1490157"

$ws.Range("B114").Value = "HP LaserJet Enterprise M607dn"
$ws.Range("C114").Value = "https://store.hp.com//us/en/pdp/hp-laserjet-enterprise-m607dn"
$ws.Range("D114").Value = $q + "949.99"
$ws.Range("E114").Value = $q + "1695184"

# --- rows 117-118 ----------------------------------------------------------
$ws.Range("B117").Value = "HP LaserJet Pro MFP M227fdw"
$ws.Range("C117").Value = "https://store.hp.com//us/en/pdp/hp-laserjet-pro-mfp-m227fdw"
$ws.Range("D117").Value = $q + "269.99"
$ws.Range("E117").Value = $q + "1490158"

$ws.Range("B118").Value = "HP LaserJet Enterprise M608x"
$ws.Range("C118").Value = "https://store.hp.com//us/en/pdp/hp-laserjet-enterprise-m608x"
$ws.Range("D118").Value = $q + "1,599.99"
$ws.Range("E118").Value = $q + "1695187"
